# Update "想去人数" (want-to-go count) figures that changed between
# the previous gh-pages data snapshot and the newly generated one.

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" ---------------------------------------------------
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 15247
$ws1.Range("F6").Value = 643
$ws1.Range("F7").Value = 1625
$ws1.Range("F8").Value = 154

# --- Sheet "全部类型" -------------------------------------------------
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 15247
$ws4.Range("F8").Value = 643
$ws4.Range("F9").Value = 1625
$ws4.Range("F11").Value = 154
